# Rename the worksheet from "Planilha1" to "Teste"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Teste"

# Remove the old header row (A1:C1 - "id", "titulo", "author") and the
# numeric id value in A2, leaving the book data (title/author columns)
# in place on rows 2 and 3.
$ws.Range("A1:C1").ClearContents()
$ws.Range("A2").ClearContents()

# Match the saved selection state (active cell C1).
$ws.Range("C1").Select()
